# Update the crop temperature table summary workbook:
# - rename the "Tbase(C)" / "Tbase_max (C)" headers to plain "Tbase" / "Tbase_max"
#   (units dropped from the header labels since the values are no longer
#   converted into different units by the processing scripts)
# - move the small colour-code legend (Black/Ecocrop, Red/Revised ...) that
#   lived in Q29:R30 over to S29:T30 so it no longer collides with the table
# - refresh the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the Tbase headers, dropping the unit suffix.
$ws.Range("Q1").Value = "Tbase"
$ws.Range("R1").Value = "Tbase_max"

# Move the legend block two columns to the right (Q->S, R->T). Copy the
# formatting only so the destination keeps the same cell style, then fill in
# the literal text, and finally remove the now-unused source cells entirely
# (not merely clear their contents) so they disappear from the sheet.
$ws.Range("Q29:R30").Copy() | Out-Null
$ws.Range("S29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("S29").Value = "Black "
$ws.Range("T29").Value = "Ecocrop"
$ws.Range("S30").Value = "Red"
$ws.Range("T30").Value = "Revised based on the new literature search"

$ws.Range("Q29:R30").Clear() | Out-Null

# Move the keyboard/selection focus like the author's session ended up.
$ws.Range("E1").Select() | Out-Null
